# Auto-generated edit script: refreshes market-price columns (H-N)
# across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets to match the
# latest scheduled price-feed snapshot.

$wb = $excel.ActiveWorkbook


$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H3").Value = 0
$ws.Range("J3").Value = 0
$ws.Range("L3").Value = 0
$ws.Range("N3").ClearContents()
$ws.Range("H62").Value = 18507.2
$ws.Range("I62").Value = 16659.125
$ws.Range("K62").Value = 16659.125
$ws.Range("M62").Value = -16035.125
$ws.Range("H65").Value = 18507.2
$ws.Range("I65").Value = 16659.125
$ws.Range("K65").Value = 83295.625
$ws.Range("M65").Value = -80175.625
$ws.Range("H102").Value = 0
$ws.Range("J102").Value = 0
$ws.Range("L102").Value = 0
$ws.Range("N102").ClearContents()
$ws.Range("H116").Value = 3579.8
$ws.Range("J116").Value = 3579.8
$ws.Range("L116").Value = 3579.8
$ws.Range("N116").Value = -10463.8
$ws.Range("H126").Value = 0
$ws.Range("J126").Value = 0
$ws.Range("L126").Value = 0
$ws.Range("N126").ClearContents()
$ws.Range("H138").Value = 28573572
$ws.Range("J138").Value = 43480212
$ws.Range("L138").Value = 130440636
$ws.Range("N138").Value = -130450916

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H21").Value = 10000
$ws.Range("I21").Value = 10000
$ws.Range("K21").Value = 10000
$ws.Range("M21").Value = -9626
$ws.Range("H32").Value = 4184.5806
$ws.Range("I32").Value = 4122.5347
$ws.Range("K32").Value = 4122.5347
$ws.Range("M32").Value = -3835.5347
$ws.Range("H48").Value = 0
$ws.Range("J48").Value = 0
$ws.Range("L48").Value = 0
$ws.Range("N48").ClearContents()
$ws.Range("H63").Value = 3313.889
$ws.Range("I63").Value = 2692.875
$ws.Range("J63").Value = 8282
$ws.Range("K63").Value = 2692.875
$ws.Range("L63").Value = 8282
$ws.Range("M63").Value = -2006.875
$ws.Range("N63").Value = -9654
$ws.Range("H66").Value = 3313.889
$ws.Range("I66").Value = 2692.875
$ws.Range("J66").Value = 8282
$ws.Range("K66").Value = 13464.375
$ws.Range("L66").Value = 41410
$ws.Range("M66").Value = -10032.375
$ws.Range("N66").Value = -48274
$ws.Range("H74").Value = 1889.1111
$ws.Range("I74").Value = 1904.0769
$ws.Range("K74").Value = 1904.0769
$ws.Range("M74").Value = -1030.0769
$ws.Range("H77").Value = 1889.1111
$ws.Range("I77").Value = 1904.0769
$ws.Range("K77").Value = 9520.3845
$ws.Range("M77").Value = -5152.3845
$ws.Range("H132").Value = 4585.75
$ws.Range("I132").Value = 4329.085
$ws.Range("K132").Value = 12987.255
$ws.Range("M132").Value = -10457.255

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 3100
$ws.Range("I22").Value = 3100
$ws.Range("K22").Value = 3100
$ws.Range("M22").Value = -2927
$ws.Range("H57").Value = 199997
$ws.Range("J57").Value = 199997
$ws.Range("L57").Value = 199997
$ws.Range("N57").Value = -201437
$ws.Range("H98").Value = 0
$ws.Range("J98").Value = 0
$ws.Range("L98").Value = 0
$ws.Range("N98").ClearContents()
$ws.Range("H134").Value = 2392.5173
$ws.Range("I134").Value = 2370.8215
$ws.Range("K134").Value = 7112.4645
$ws.Range("M134").Value = -4577.4645
$ws.Range("H136").Value = 199997
$ws.Range("J136").Value = 199997
$ws.Range("L136").Value = 199997
$ws.Range("N136").Value = -210197

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 6910.6113
$ws.Range("I58").Value = 2824.5
$ws.Range("J58").Value = 10179.5
$ws.Range("K58").Value = 2824.5
$ws.Range("L58").Value = 10179.5
$ws.Range("M58").Value = -2621.5
$ws.Range("N58").Value = -10585.5
$ws.Range("H88").Value = 13051.833
$ws.Range("J88").Value = 13051.833
$ws.Range("L88").Value = 13051.833
$ws.Range("N88").Value = -13863.833
$ws.Range("H91").Value = 13051.833
$ws.Range("J91").Value = 13051.833
$ws.Range("L91").Value = 13051.833
$ws.Range("N91").Value = -15859.833
$ws.Range("H136").Value = 6910.6113
$ws.Range("I136").Value = 2824.5
$ws.Range("J136").Value = 10179.5
$ws.Range("K136").Value = 8473.5
$ws.Range("L136").Value = 30538.5
$ws.Range("M136").Value = -5923.5
$ws.Range("N136").Value = -35638.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1029.8889
$ws.Range("J5").Value = 1145.8334
$ws.Range("L5").Value = 3437.5002
$ws.Range("N5").Value = -3661.5002
$ws.Range("H40").Value = 267.63635
$ws.Range("I40").Value = 310.5
$ws.Range("J40").Value = 216.2
$ws.Range("K40").Value = 1242
$ws.Range("L40").Value = 864.8
$ws.Range("M40").Value = -1173
$ws.Range("N40").Value = -1002.8
$ws.Range("H132").Value = 35715340
$ws.Range("I132").Value = 45455376
$ws.Range("J132").Value = 1866
$ws.Range("K132").Value = 409098384
$ws.Range("L132").Value = 16794
$ws.Range("M132").Value = -409095854
$ws.Range("N132").Value = -21854
$ws.Range("H135").Value = 1029.8889
$ws.Range("J135").Value = 1145.8334
$ws.Range("L135").Value = 10312.5006
$ws.Range("N135").Value = -15382.5006
$ws.Range("H137").Value = 6631.231
$ws.Range("I137").Value = 1354
$ws.Range("J137").Value = 24222
$ws.Range("K137").Value = 4062
$ws.Range("L137").Value = 72666
$ws.Range("M137").Value = 1038
$ws.Range("N137").Value = -82866

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 10292.3
$ws.Range("I2").Value = 314.75
$ws.Range("J2").Value = 16944
$ws.Range("K2").Value = 314.75
$ws.Range("L2").Value = 16944
$ws.Range("M2").Value = -201.75
$ws.Range("N2").Value = -17170
$ws.Range("H20").Value = 17955.5
$ws.Range("J20").Value = 17955.5
$ws.Range("L20").Value = 17955.5
$ws.Range("N20").Value = -18445.5
$ws.Range("H31").Value = 9775
$ws.Range("I31").Value = 9775
$ws.Range("K31").Value = 9775
$ws.Range("M31").Value = -9483
$ws.Range("H37").Value = 9775
$ws.Range("I37").Value = 9775
$ws.Range("K37").Value = 9775
$ws.Range("M37").Value = -9498
$ws.Range("H80").Value = 4492.095
$ws.Range("I80").Value = 3937.077
$ws.Range("K80").Value = 3937.077
$ws.Range("M80").Value = -2939.077
$ws.Range("H83").Value = 4492.095
$ws.Range("I83").Value = 3937.077
$ws.Range("K83").Value = 19685.385
$ws.Range("M83").Value = -14693.385
$ws.Range("H113").Value = 3761.7058
$ws.Range("I113").Value = 3301.1428
$ws.Range("K113").Value = 3301.1428
$ws.Range("M113").Value = -1131.1428

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H5").Value = 11011
$ws.Range("J5").Value = 11011
$ws.Range("L5").Value = 11011
$ws.Range("N5").Value = -11237
$ws.Range("H22").Value = 3206.125
$ws.Range("J22").Value = 3837.5
$ws.Range("L22").Value = 3837.5
$ws.Range("N22").Value = -4427.5
$ws.Range("H27").Value = 3206.125
$ws.Range("J27").Value = 3837.5
$ws.Range("L27").Value = 3837.5
$ws.Range("N27").Value = -4051.5
$ws.Range("H55").Value = 938
$ws.Range("I55").Value = 958.7273
$ws.Range("J55").Value = 900
$ws.Range("K55").Value = 958.7273
$ws.Range("L55").Value = 900
$ws.Range("M55").Value = -785.7273
$ws.Range("N55").Value = -1246
$ws.Range("H68").Value = 16943.5
$ws.Range("J68").Value = 16616.334
$ws.Range("L68").Value = 16616.334
$ws.Range("N68").Value = -18114.334
$ws.Range("H71").Value = 16943.5
$ws.Range("J71").Value = 16616.334
$ws.Range("L71").Value = 83081.67
$ws.Range("N71").Value = -90569.67
$ws.Range("H82").Value = 8313.931
$ws.Range("I82").Value = 10139.625
$ws.Range("K82").Value = 10139.625
$ws.Range("M82").Value = -9778.625
$ws.Range("H85").Value = 8313.931
$ws.Range("I85").Value = 10139.625
$ws.Range("K85").Value = 10139.625
$ws.Range("M85").Value = -8891.625
$ws.Range("H132").Value = 16320.333
$ws.Range("I132").Value = 26154.666
$ws.Range("J132").Value = 6486
$ws.Range("K132").Value = 78463.99800000001
$ws.Range("L132").Value = 19458
$ws.Range("M132").Value = -75933.99800000001
$ws.Range("N132").Value = -24518

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 4100
$ws.Range("I81").Value = 2625.25
$ws.Range("K81").Value = 5250.5
$ws.Range("M81").Value = -4189.5
$ws.Range("H84").Value = 4100
$ws.Range("I84").Value = 2625.25
$ws.Range("K84").Value = 26252.5
$ws.Range("M84").Value = -20948.5
$ws.Range("H106").Value = 44980.2
$ws.Range("J106").Value = 44980.2
$ws.Range("L106").Value = 44980.2
$ws.Range("N106").Value = -47504.2
$ws.Range("H107").Value = 2166.923
$ws.Range("I107").Value = 982.7778
$ws.Range("J107").Value = 4831.25
$ws.Range("K107").Value = 2948.3334
$ws.Range("L107").Value = 14493.75
$ws.Range("M107").Value = -1028.3334
$ws.Range("N107").Value = -18333.75
$ws.Range("H113").Value = 569.8333
$ws.Range("I113").Value = 579.55
$ws.Range("J113").Value = 521.25
$ws.Range("K113").Value = 1738.65
$ws.Range("L113").Value = 1563.75
$ws.Range("M113").Value = 431.3500000000001
$ws.Range("N113").Value = -5903.75
